$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: splits a run whose text is "Operador <rest>" into two runs with
# identical formatting: "Ator " and "<rest>". We cannot just assign a new
# value to Range.Text because the engine silently re-merges two adjacent
# runs that end up with identical rPr. Toggling a character property
# (Bold on/off) after each piece is written is what keeps the run boundary
# in the saved package, so we do that for both halves.
# ---------------------------------------------------------------------------
function Split-OperadorRun([string]$oldText, [string]$restText) {
    $rng = $d.Content
    $rng.Find.Execute($oldText) | Out-Null
    $rng.Collapse(1)
    $rng.InsertAfter("Ator ")

    $atorRange = $d.Range($rng.Start, $rng.Start + 5)
    $atorRange.Bold = 1
    $atorRange.Bold = 0

    # InsertAfter only inserts "Ator " - the untouched original text ("Operador
    # <rest>") is still sitting right after it, so the range to overwrite spans
    # the *original* text's full length, not the trimmed replacement's length.
    $restRange = $d.Range($atorRange.End, $atorRange.End + $oldText.Length)
    $restRange.Text = $restText
    $restRange.Bold = 1
    $restRange.Bold = 0
}

Split-OperadorRun "Operador clica no botão gerenciar máquinas." "clica no botão gerenciar máquinas."
Split-OperadorRun "Operador seleciona cliente desejado" "seleciona cliente desejado"
Split-OperadorRun "Operador seleciona a máquina desejada." "seleciona a máquina desejada."
Split-OperadorRun "Operador clica no botão “Abrir máquina”" "clica no botão “Abrir máquina”"
Split-OperadorRun "Operador seleciona o ícone com o sinal de adição para " "seleciona o ícone com o sinal de adição para "
Split-OperadorRun "Operador preenche as informações necessárias e clica no botão adicionar" "preenche as informações necessárias e clica no botão adicionar"

# ---------------------------------------------------------------------------
# "Sistema exibe uma lista com as" / " máquinas vinculadas ao cliente." is
# merged back into a single run. The two original runs are separated by the
# hidden "_GoBack" bookmark, so we drive the merge through Find & Replace
# (rather than a manual Range.Text assignment) because that is what makes
# the engine correctly discard the bookmark the replaced span swallows.
# ---------------------------------------------------------------------------
$mergeRng = $d.Content
$mergeRng.Find.Execute(
    "Sistema exibe uma lista com as máquinas vinculadas ao cliente.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Sistema exibe uma lista com as máquinas vinculadas ao cliente.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Last occurrence also carries the relocated "_GoBack" bookmark, now sitting
# right after the new "Ator " run.
# ---------------------------------------------------------------------------
Split-OperadorRun "Operador clica no botão “Abrir máquina”." "clica no botão “Abrir máquina”."

$bmRng = $d.Content
$bmRng.Find.Execute("clica no botão “Abrir máquina”.") | Out-Null
$bmPoint = $d.Range($bmRng.Start, $bmRng.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

Write-Output "Done"
